$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 3.405752891157269
$ws.Cells.Item(2, 4).Value = 4.363120792079042
$ws.Cells.Item(2, 5).Value = 11.37321447006714
$ws.Cells.Item(2, 6).Value = 22.72240070235265
$ws.Cells.Item(2, 7).Value = 3.614682361324854
$ws.Cells.Item(2, 9).Value = 18.99795226744325
$ws.Cells.Item(2, 11).Value = 11.64681834076186
$ws.Cells.Item(2, 13).Value = 14.67633183042082
$ws.Cells.Item(2, 14).Value = 16.93502142271931
$ws.Cells.Item(2, 15).Value = 20.21510641877277

$ws.Cells.Item(3, 3).Value = 3.28883344543926
$ws.Cells.Item(3, 4).Value = 4.320569732506225
$ws.Cells.Item(3, 5).Value = 11.26614973879036
$ws.Cells.Item(3, 6).Value = 22.68770855210153
$ws.Cells.Item(3, 7).Value = 3.616738160184558
$ws.Cells.Item(3, 9).Value = 19.01050877108378
$ws.Cells.Item(3, 11).Value = 11.18046709010335
$ws.Cells.Item(3, 13).Value = 14.39765365074639
$ws.Cells.Item(3, 14).Value = 16.98145657472332
$ws.Cells.Item(3, 15).Value = 20.25071813394637

$ws.Cells.Item(4, 3).Value = 3.214178232495343
$ws.Cells.Item(4, 4).Value = 4.293895205800359
$ws.Cells.Item(4, 5).Value = 11.20452287219994
$ws.Cells.Item(4, 6).Value = 22.67387577960579
$ws.Cells.Item(4, 7).Value = 3.618067423836294
$ws.Cells.Item(4, 9).Value = 19.02370649679935
$ws.Cells.Item(4, 11).Value = 10.88501284462108
$ws.Cells.Item(4, 13).Value = 14.22702968079303
$ws.Cells.Item(4, 14).Value = 17.01175660176978
$ws.Cells.Item(4, 15).Value = 20.27831474234597

$ws.Cells.Item(5, 3).Value = 3.183061407480082
$ws.Cells.Item(5, 4).Value = 4.282892555894665
$ws.Cells.Item(5, 5).Value = 11.1804700447145
$ws.Cells.Item(5, 6).Value = 22.67011955852888
$ws.Cells.Item(5, 7).Value = 3.618626010040281
$ws.Cells.Item(5, 9).Value = 19.03046198997574
$ws.Cells.Item(5, 11).Value = 10.76249491872392
$ws.Cells.Item(5, 13).Value = 14.15771812157703
$ws.Cells.Item(5, 14).Value = 17.02455470398252
$ws.Cells.Item(5, 15).Value = 20.29099680558446

$ws.Cells.Item(6, 3).Value = 3.177853386292623
$ws.Cells.Item(6, 4).Value = 4.28105772220007
$ws.Cells.Item(6, 5).Value = 11.17654088038108
$ws.Cells.Item(6, 6).Value = 22.66960947359577
$ws.Cells.Item(6, 7).Value = 3.618719785192143
$ws.Cells.Item(6, 9).Value = 19.03166682603455
$ws.Cells.Item(6, 11).Value = 10.74202846931732
$ws.Cells.Item(6, 13).Value = 14.14622508109598
$ws.Cells.Item(6, 14).Value = 17.02670705848136
$ws.Cells.Item(6, 15).Value = 20.29318924338982

$ws.Cells.Item(7, 3).Value = 3.213761353782027
$ws.Cells.Item(7, 4).Value = 4.293747350000788
$ws.Cells.Item(7, 5).Value = 11.20419415953207
$ws.Cells.Item(7, 6).Value = 22.67381750474171
$ws.Cells.Item(7, 7).Value = 3.618074888623969
$ws.Cells.Item(7, 9).Value = 19.02379203124739
$ws.Cells.Item(7, 11).Value = 10.8833688559996
$ws.Cells.Item(7, 13).Value = 14.22609390470794
$ws.Cells.Item(7, 14).Value = 17.01192737591562
$ws.Cells.Item(7, 15).Value = 20.27847996876359

$ws.Cells.Item(8, 3).Value = 3.366048942981504
$ws.Cells.Item(8, 4).Value = 4.348566039128088
$ws.Cells.Item(8, 5).Value = 11.3354621339342
$ws.Cells.Item(8, 6).Value = 22.70889107660476
$ws.Cells.Item(8, 7).Value = 3.615377328597137
$ws.Cells.Item(8, 9).Value = 19.0011413507584
$ws.Cells.Item(8, 11).Value = 11.48800960831236
$ws.Cells.Item(8, 13).Value = 14.58020056589688
$ws.Cells.Item(8, 14).Value = 16.95066151911655
$ws.Cells.Item(8, 15).Value = 20.22619315186271

$ws.Cells.Item(9, 3).Value = 3.640924161902158
$ws.Cells.Item(9, 4).Value = 4.45149267734341
$ws.Cells.Item(9, 5).Value = 11.62419618043178
$ws.Cells.Item(9, 6).Value = 22.83671199944093
$ws.Cells.Item(9, 7).Value = 3.61061652728784
$ws.Cells.Item(9, 9).Value = 19.00035203348218
$ws.Cells.Item(9, 11).Value = 12.59478571171978
$ws.Cells.Item(9, 13).Value = 15.27428087514819
$ws.Cells.Item(9, 14).Value = 16.84467491170499
$ws.Cells.Item(9, 15).Value = 20.1693222717478

$ws.Cells.Item(10, 3).Value = 3.827215889791427
$ws.Cells.Item(10, 4).Value = 4.524033182443152
$ws.Cells.Item(10, 5).Value = 11.85351900416002
$ws.Cells.Item(10, 6).Value = 22.96618839168923
$ws.Cells.Item(10, 7).Value = 3.607437867507929
$ws.Cells.Item(10, 9).Value = 19.02643063086854
$ws.Cells.Item(10, 11).Value = 13.35225490968668
$ws.Cells.Item(10, 13).Value = 15.77861400992806
$ws.Cells.Item(10, 14).Value = 16.77538621346867
$ws.Cells.Item(10, 15).Value = 20.15560036541045

$ws.Cells.Item(11, 3).Value = 3.908365020065501
$ws.Cells.Item(11, 4).Value = 4.556305290791151
$ws.Cells.Item(11, 5).Value = 11.961135119094
$ws.Cells.Item(11, 6).Value = 23.03267854186665
$ws.Cells.Item(11, 7).Value = 3.606060367171962
$ws.Cells.Item(11, 9).Value = 19.04407730754018
$ws.Cells.Item(11, 11).Value = 13.68349126087759
$ws.Cells.Item(11, 13).Value = 16.00575156833993
$ws.Cells.Item(11, 14).Value = 16.74571777656648
$ws.Cells.Item(11, 15).Value = 20.15548498789976

$ws.Cells.Item(12, 3).Value = 3.938563121533847
$ws.Cells.Item(12, 4).Value = 4.568416476342971
$ws.Cells.Item(12, 5).Value = 12.00232174528756
$ws.Cells.Item(12, 6).Value = 23.05893353664048
$ws.Cells.Item(12, 7).Value = 3.60554853633377
$ws.Cells.Item(12, 9).Value = 19.05158917499997
$ws.Cells.Item(12, 11).Value = 13.80691727569821
$ws.Cells.Item(12, 13).Value = 16.09134951303901
$ws.Cells.Item(12, 14).Value = 16.73474860293804
$ws.Cells.Item(12, 15).Value = 20.15632387292947

$ws.Cells.Item(13, 3).Value = 3.932083256702676
$ws.Cells.Item(13, 4).Value = 4.565813067161997
$ws.Cells.Item(13, 5).Value = 11.99343276579177
$ws.Cells.Item(13, 6).Value = 23.05323142947287
$ws.Cells.Item(13, 7).Value = 3.605658333254678
$ws.Cells.Item(13, 9).Value = 19.04993451489716
$ws.Cells.Item(13, 11).Value = 13.78042575518468
$ws.Cells.Item(13, 13).Value = 16.0729341922999
$ws.Cells.Item(13, 14).Value = 16.73709920654881
$ws.Cells.Item(13, 15).Value = 20.15610393454313

$ws.Cells.Item(14, 3).Value = 3.91086015907316
$ws.Cells.Item(14, 4).Value = 4.557303909388506
$ws.Cells.Item(14, 5).Value = 11.96451508217674
$ws.Cells.Item(14, 6).Value = 23.03481707155681
$ws.Cells.Item(14, 7).Value = 3.60601806244599
$ws.Cells.Item(14, 9).Value = 19.04467870588416
$ws.Cells.Item(14, 11).Value = 13.69368625413019
$ws.Cells.Item(14, 13).Value = 16.01280253491909
$ws.Cells.Item(14, 14).Value = 16.74481001599506
$ws.Cells.Item(14, 15).Value = 20.15553630913365

$ws.Cells.Item(15, 3).Value = 3.897790816731351
$ws.Cells.Item(15, 4).Value = 4.552077389592796
$ws.Cells.Item(15, 5).Value = 11.94685760447393
$ws.Cells.Item(15, 6).Value = 23.02367747642891
$ws.Cells.Item(15, 7).Value = 3.606239681561247
$ws.Cells.Item(15, 9).Value = 19.0415673090023
$ws.Cells.Item(15, 11).Value = 13.64029217810919
$ws.Cells.Item(15, 13).Value = 15.97591377286345
$ws.Cells.Item(15, 14).Value = 16.74956768351159
$ws.Cells.Item(15, 15).Value = 20.15530359280968

$ws.Cells.Item(16, 3).Value = 3.82183898421283
$ws.Cells.Item(16, 4).Value = 4.521909063353882
$ws.Cells.Item(16, 5).Value = 11.84654902807544
$ws.Cells.Item(16, 6).Value = 22.96199461895393
$ws.Cells.Item(16, 7).Value = 3.607529264277773
$ws.Cells.Item(16, 9).Value = 19.02539359045253
$ws.Cells.Item(16, 11).Value = 13.33033160597131
$ws.Cells.Item(16, 13).Value = 15.76371663084941
$ws.Cells.Item(16, 14).Value = 16.77736231736156
$ws.Cells.Item(16, 15).Value = 20.1557313420005

$ws.Cells.Item(17, 3).Value = 3.774313254794385
$ws.Cells.Item(17, 4).Value = 4.503211964443754
$ws.Cells.Item(17, 5).Value = 11.78582952081091
$ws.Cells.Item(17, 6).Value = 22.92608836935137
$ws.Cells.Item(17, 7).Value = 3.608337887328032
$ws.Cells.Item(17, 9).Value = 19.01695155119998
$ws.Cells.Item(17, 11).Value = 13.13669756235144
$ws.Cells.Item(17, 13).Value = 15.63289283766983
$ws.Cells.Item(17, 14).Value = 16.79488716222658
$ws.Cells.Item(17, 15).Value = 20.15756431129196

$ws.Cells.Item(18, 3).Value = 3.746639968735194
$ws.Cells.Item(18, 4).Value = 4.492389757532822
$ws.Cells.Item(18, 5).Value = 11.75121791475047
$ws.Cells.Item(18, 6).Value = 22.90615153827176
$ws.Cells.Item(18, 7).Value = 3.60880943528282
$ws.Cells.Item(18, 9).Value = 19.0126404660126
$ws.Cells.Item(18, 11).Value = 13.02407311251049
$ws.Cells.Item(18, 13).Value = 15.55743609779812
$ws.Cells.Item(18, 14).Value = 16.80514128064501
$ws.Cells.Item(18, 15).Value = 20.15919522320148

$ws.Cells.Item(19, 3).Value = 3.737212727064418
$ws.Cells.Item(19, 4).Value = 4.488713998387944
$ws.Cells.Item(19, 5).Value = 11.7395538985816
$ws.Cells.Item(19, 6).Value = 22.8995245781684
$ws.Cells.Item(19, 7).Value = 3.608970202527476
$ws.Cells.Item(19, 9).Value = 19.01127438502986
$ws.Cells.Item(19, 11).Value = 12.98572838641502
$ws.Cells.Item(19, 13).Value = 15.53185428995272
$ws.Cells.Item(19, 14).Value = 16.8086431023489
$ws.Cells.Item(19, 15).Value = 20.1598464002881

$ws.Cells.Item(20, 3).Value = 3.779407514145096
$ws.Cells.Item(20, 4).Value = 4.505209390631258
$ws.Cells.Item(20, 5).Value = 11.79226113961785
$ws.Cells.Item(20, 6).Value = 22.92983670162479
$ws.Cells.Item(20, 7).Value = 3.608251140900876
$ws.Cells.Item(20, 9).Value = 19.01779387432939
$ws.Cells.Item(20, 11).Value = 13.15744041435241
$ws.Cells.Item(20, 13).Value = 15.64684165769143
$ws.Cells.Item(20, 14).Value = 16.79300357850773
$ws.Cells.Item(20, 15).Value = 20.1573094952166

$ws.Cells.Item(21, 3).Value = 3.917108420923444
$ws.Cells.Item(21, 4).Value = 4.559806268128638
$ws.Cells.Item(21, 5).Value = 11.97299742450691
$ws.Cells.Item(21, 6).Value = 23.04019672446057
$ws.Cells.Item(21, 7).Value = 3.605912135725485
$ws.Cells.Item(21, 9).Value = 19.04619997549429
$ws.Cells.Item(21, 11).Value = 13.71921882683321
$ws.Cells.Item(21, 13).Value = 16.03047655685742
$ws.Cells.Item(21, 14).Value = 16.74253796017729
$ws.Cells.Item(21, 15).Value = 20.15567907298936

$ws.Cells.Item(22, 3).Value = 4.004001730630779
$ws.Cells.Item(22, 4).Value = 4.594847645697819
$ws.Cells.Item(22, 5).Value = 12.09363436769966
$ws.Cells.Item(22, 6).Value = 23.11859042816574
$ws.Cells.Item(22, 7).Value = 3.604440551761814
$ws.Cells.Item(22, 9).Value = 19.06959799987243
$ws.Cells.Item(22, 11).Value = 14.07464891909435
$ws.Cells.Item(22, 13).Value = 16.27875404662951
$ws.Cells.Item(22, 14).Value = 16.7111037336064
$ws.Cells.Item(22, 15).Value = 20.15975801053102

$ws.Cells.Item(23, 3).Value = 3.957913203140938
$ws.Cells.Item(23, 4).Value = 4.576205663803183
$ws.Cells.Item(23, 5).Value = 12.02903114642216
$ws.Cells.Item(23, 6).Value = 23.0761822346877
$ws.Cells.Item(23, 7).Value = 3.605220756345312
$ws.Cells.Item(23, 9).Value = 19.05666878767752
$ws.Cells.Item(23, 11).Value = 13.88604776163826
$ws.Cells.Item(23, 13).Value = 16.14649486180766
$ws.Cells.Item(23, 14).Value = 16.72773932459424
$ws.Cells.Item(23, 15).Value = 20.15710996452474

$ws.Cells.Item(24, 3).Value = 3.777105488924119
$ws.Cells.Item(24, 4).Value = 4.504306581225334
$ws.Cells.Item(24, 5).Value = 11.78935247713716
$ws.Cells.Item(24, 6).Value = 22.92813987985518
$ws.Cells.Item(24, 7).Value = 3.608290338210973
$ws.Cells.Item(24, 9).Value = 19.01741137045833
$ws.Cells.Item(24, 11).Value = 13.14806661854417
$ws.Cells.Item(24, 13).Value = 15.64053615084649
$ws.Cells.Item(24, 14).Value = 16.79385458952776
$ws.Cells.Item(24, 15).Value = 20.1574228999811

$ws.Cells.Item(25, 3).Value = 3.569237791023917
$ws.Cells.Item(25, 4).Value = 4.424168394244588
$ws.Cells.Item(25, 5).Value = 11.54292204000557
$ws.Cells.Item(25, 6).Value = 22.79584794211297
$ws.Cells.Item(25, 7).Value = 3.611848164434801
$ws.Cells.Item(25, 9).Value = 18.99588487817445
$ws.Cells.Item(25, 11).Value = 12.30467086564377
$ws.Cells.Item(25, 13).Value = 15.27428087514819
$ws.Cells.Item(25, 14).Value = 16.84467491170499
$ws.Cells.Item(25, 15).Value = 20.17979303731479
